$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 1630
$ws.Range("I127").Value = 776.6667
$ws.Range("K127").Value = 2330.0001
$ws.Range("M127").Value = 2629.9999
$ws.Range("H129").Value = 1725.909
$ws.Range("I129").Value = 515
$ws.Range("J129").Value = 2180
$ws.Range("K129").Value = 1545
$ws.Range("L129").Value = 6540
$ws.Range("M129").Value = 3455
$ws.Range("N129").Value = -16540
$ws.Range("H141").Value = 2496.7368
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3092.7222
$ws.Range("I2").Value = 1200.8889
$ws.Range("J2").Value = 4984.5557
$ws.Range("K2").Value = 1200.8889
$ws.Range("L2").Value = 4984.5557
$ws.Range("M2").Value = -1087.8889
$ws.Range("N2").Value = -5210.5557
$ws.Range("H32").Value = 2318403
$ws.Range("I32").Value = 2406605
$ws.Range("K32").Value = 2406605
$ws.Range("M32").Value = -2406318
$ws.Range("H74").Value = 32351.258
$ws.Range("I74").Value = 39627.445
$ws.Range("J74").Value = 7794.125
$ws.Range("K74").Value = 39627.445
$ws.Range("L74").Value = 7794.125
$ws.Range("M74").Value = -38753.445
$ws.Range("N74").Value = -9542.125
$ws.Range("H77").Value = 32351.258
$ws.Range("I77").Value = 39627.445
$ws.Range("J77").Value = 7794.125
$ws.Range("K77").Value = 198137.225
$ws.Range("L77").Value = 38970.625
$ws.Range("M77").Value = -193769.225
$ws.Range("N77").Value = -47706.625
$ws.Range("H116").Value = 3092.7222
$ws.Range("I116").Value = 1200.8889
$ws.Range("J116").Value = 4984.5557
$ws.Range("K116").Value = 1200.8889
$ws.Range("L116").Value = 4984.5557
$ws.Range("M116").Value = 1093.1111
$ws.Range("N116").Value = -9572.555700000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3092.7222
$ws.Range("I3").Value = 1200.8889
$ws.Range("J3").Value = 4984.5557
$ws.Range("K3").Value = 1200.8889
$ws.Range("L3").Value = 4984.5557
$ws.Range("M3").Value = -1086.8889
$ws.Range("N3").Value = -5212.5557
$ws.Range("H134").Value = 8966.883
$ws.Range("I134").Value = 3555.375
$ws.Range("K134").Value = 10666.125
$ws.Range("M134").Value = -8131.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 184.4
$ws.Range("I7").Value = 62.75
$ws.Range("J7").Value = 366.875
$ws.Range("K7").Value = 62.75
$ws.Range("L7").Value = 366.875
$ws.Range("M7").Value = 50.25
$ws.Range("N7").Value = -592.875
$ws.Range("H16").Value = 4242.839
$ws.Range("I16").Value = 3780.05
$ws.Range("K16").Value = 3780.05
$ws.Range("M16").Value = -3493.05
$ws.Range("H22").Value = 325.8
$ws.Range("I22").Value = 295.33334
$ws.Range("K22").Value = 295.33334
$ws.Range("M22").Value = 54.66665999999998
$ws.Range("H31").Value = 11922.5
$ws.Range("I31").Value = 3700.75
$ws.Range("J31").Value = 18499.9
$ws.Range("K31").Value = 3700.75
$ws.Range("L31").Value = 18499.9
$ws.Range("M31").Value = -3405.75
$ws.Range("N31").Value = -19089.9
$ws.Range("H34").Value = 11922.5
$ws.Range("I34").Value = 3700.75
$ws.Range("J34").Value = 18499.9
$ws.Range("K34").Value = 3700.75
$ws.Range("L34").Value = 18499.9
$ws.Range("M34").Value = -3498.75
$ws.Range("N34").Value = -18903.9
$ws.Range("H113").Value = 4242.839
$ws.Range("I113").Value = 3780.05
$ws.Range("K113").Value = 3780.05
$ws.Range("M113").Value = -1610.05
$ws.Range("H122").Value = 3350
$ws.Range("I122").Value = 3260
$ws.Range("K122").Value = 9780
$ws.Range("M122").Value = -7330
$ws.Range("H134").Value = 5648.6294
$ws.Range("I134").Value = 2177.2
$ws.Range("J134").Value = 15567
$ws.Range("K134").Value = 6531.599999999999
$ws.Range("L134").Value = 46701
$ws.Range("M134").Value = -3996.599999999999
$ws.Range("N134").Value = -51771

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1786.0541
$ws.Range("I5").Value = 1584.7084
$ws.Range("J5").Value = 2157.7693
$ws.Range("K5").Value = 4754.1252
$ws.Range("L5").Value = 6473.3079
$ws.Range("M5").Value = -4642.1252
$ws.Range("N5").Value = -6697.3079
$ws.Range("H40").Value = 792.9231
$ws.Range("I40").Value = 13.222222
$ws.Range("K40").Value = 52.888888
$ws.Range("M40").Value = 16.111112
$ws.Range("H41").Value = 1162
$ws.Range("I41").Value = 1000
$ws.Range("J41").Value = 1208.2858
$ws.Range("K41").Value = 3000
$ws.Range("L41").Value = 3624.8574
$ws.Range("M41").Value = -2662
$ws.Range("N41").Value = -4300.857400000001
$ws.Range("H51").Value = 887
$ws.Range("J51").Value = 965.7
$ws.Range("L51").Value = 2897.1
$ws.Range("N51").Value = -3817.1
$ws.Range("H68").Value = 2226.423
$ws.Range("J68").Value = 2423.842
$ws.Range("L68").Value = 7271.526
$ws.Range("N68").Value = -8893.526
$ws.Range("H71").Value = 2226.423
$ws.Range("J71").Value = 2423.842
$ws.Range("L71").Value = 21814.578
$ws.Range("N71").Value = -29926.578
$ws.Range("H107").Value = 426808.47
$ws.Range("J107").Value = 279251.38
$ws.Range("L107").Value = 837754.14
$ws.Range("N107").Value = -841594.14
$ws.Range("H113").Value = 1878.0555
$ws.Range("I113").Value = 1152.25
$ws.Range("J113").Value = 2458.7
$ws.Range("K113").Value = 3456.75
$ws.Range("L113").Value = 7376.099999999999
$ws.Range("M113").Value = -1286.75
$ws.Range("N113").Value = -11716.1
$ws.Range("H119").Value = 9995
$ws.Range("I119").Value = 9995
$ws.Range("K119").Value = 29985
$ws.Range("M119").Value = -25147
$ws.Range("H124").Value = 875
$ws.Range("I124").Value = 875
$ws.Range("K124").Value = 2625
$ws.Range("M124").Value = 2285
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("H129").Value = 10479804
$ws.Range("J129").Value = 33534240
$ws.Range("L129").Value = 100602720
$ws.Range("N129").Value = -100612720
$ws.Range("H131").Value = 2119.5898
$ws.Range("I131").Value = 1919.1666
$ws.Range("J131").Value = 2156.0303
$ws.Range("K131").Value = 5757.4998
$ws.Range("L131").Value = 6468.090899999999
$ws.Range("M131").Value = -717.4997999999996
$ws.Range("N131").Value = -16548.0909
$ws.Range("H135").Value = 1786.0541
$ws.Range("I135").Value = 1584.7084
$ws.Range("J135").Value = 2157.7693
$ws.Range("K135").Value = 14262.3756
$ws.Range("L135").Value = 19419.9237
$ws.Range("M135").Value = -11727.3756
$ws.Range("N135").Value = -24489.9237

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 6666863.5
$ws.Range("I2").Value = 100
$ws.Range("K2").Value = 100
$ws.Range("M2").Value = 13

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3391.4
$ws.Range("J16").Value = 3999.5
$ws.Range("L16").Value = 3999.5
$ws.Range("N16").Value = -4339.5
$ws.Range("H40").Value = 62504412
$ws.Range("I40").Value = 62504412
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 62504412
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -62504276
$ws.Range("N40").ClearContents()
$ws.Range("H46").Value = 2853.9092
$ws.Range("I46").Value = 1378.2
$ws.Range("K46").Value = 1378.2
$ws.Range("M46").Value = -1190.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H55").Value = 29000
$ws.Range("I55").Value = 48000
$ws.Range("K55").Value = 48000
$ws.Range("M55").Value = -47723
$ws.Range("H107").Value = 836.3570999999999
$ws.Range("I107").Value = 783
$ws.Range("K107").Value = 2349
$ws.Range("M107").Value = -429
